$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 200
$ws1.Range("F3").Value = 544
$ws1.Range("F9").Value = 399
$ws1.Range("F10").Value = 3448
$ws1.Range("F11").Value = 47

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 200
$ws4.Range("F4").Value = 544
$ws4.Range("F10").Value = 399
$ws4.Range("F11").Value = 3448
$ws4.Range("F12").Value = 47
